$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.703.43'
$ws.Range("E2").Value = '  -0.78%  '

# Row 3
$ws.Range("D3").Value = '1.631.40'
$ws.Range("E3").Value = '  -1.13%  '

# Row 4
$ws.Range("E4").Value = '  -0.80%  '

# Row 5
$ws.Range("D5").Value = '218.75'
$ws.Range("E5").Value = '  +0.36%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.497'
$ws.Range("E6").Value = '  -1.34%  '

# Row 7
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  -0.75%  '

# Row 8
$ws.Range("E8").Value = '  -1.54%  '

# Row 9
$ws.Range("D9").Value = '0.0618'
$ws.Range("E9").Value = '  -1.53%  '

# Row 10
$ws.Range("D10").Value = '18.89'
$ws.Range("E10").Value = '  -1.48%  '

# Row 11
$ws.Range("E11").Value = '  -0.17%  '

# Row 12
$ws.Range("D12").Value = '1.858.60'
$ws.Range("E12").Value = '  -1.09%  '

# Row 13
$ws.Range("D13").Value = '1.622.18'
$ws.Range("E13").Value = '  -1.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  -2.47%  '

# Row 15
$ws.Range("E15").Value = '  -1.60%  '

# Row 16
$ws.Range("D16").Value = '64.17'
$ws.Range("E16").Value = '  -1.21%  '

# Row 17
$ws.Range("D17").Value = '26.672.21'
$ws.Range("E17").Value = '  -0.89%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -1.78%  '

# Row 19
$ws.Range("D19").Value = '212.67'
$ws.Range("E19").Value = '  -0.48%  '

# Row 20
$ws.Range("E20").Value = '  -0.66%  '

# Row 21
$ws.Range("D21").Value = '4.32'
$ws.Range("E21").Value = '  -1.03%  '

# Row 22
$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  -1.36%  '

# Row 23
$ws.Range("E23").Value = '  -2.70%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.00'
$ws.Range("E24").Value = '  -4.18%  '

# Row 25
$ws.Range("D25").Value = '147.88'
$ws.Range("E25").Value = '  +1.20%  '

# Row 26
$ws.Range("E26").Value = '  -0.67%  '

# Row 27
$ws.Range("D27").Value = '0.118'
$ws.Range("E27").Value = '  -0.53%  '

# Row 28
$ws.Range("D28").Value = '7.02'
$ws.Range("E28").Value = '  -1.41%  '

# Row 29
$ws.Range("D29").Value = '15.55'
$ws.Range("E29").Value = '  -1.19%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0500'
$ws.Range("E30").Value = '  -3.11%  '

# Row 31
$ws.Range("E31").Value = '  +0.61%  '

# Row 32
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  +1.44%  '

# Row 33
$ws.Range("E33").Value = '  -1.14%  '

# Row 34
$ws.Range("E34").Value = '  -0.82%  '

# Row 35
$ws.Range("D35").Value = '1.254.30'
$ws.Range("E35").Value = '  -1.87%  '

# Row 36
$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("E37").Value = '  -1.48%  '

# Row 38
$ws.Range("E38").Value = '  -2.67%  '

# Row 39
$ws.Range("E39").Value = '  -0.73%  '

# Row 40
$ws.Range("D40").Value = '0.803'
$ws.Range("E40").Value = '  -3.26%  '

# Row 41
$ws.Range("E41").Value = '  -1.83%  '

# Row 42
$ws.Range("E42").Value = '  -2.04%  '

# Row 43
$ws.Range("D43").Value = '1.768.60'
$ws.Range("E43").Value = '  -1.71%  '

# Row 44
$ws.Range("D44").Value = '2.12'
$ws.Range("E44").Value = '  -5.79%  '

# Row 45
$ws.Range("D45").Value = '91.81'
$ws.Range("E45").Value = '  +0.07%  '

# Row 46
$ws.Range("D46").Value = '59.54'
$ws.Range("E46").Value = '  +1.30%  '

# Row 47
$ws.Range("D47").Value = '1.57'
$ws.Range("E47").Value = '  -2.47%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0102'
$ws.Range("E48").Value = '  -2.32%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0516'
$ws.Range("E49").Value = '  -0.97%  '

# Row 50
$ws.Range("E50").Value = '  -0.69%  '

# Row 51
$ws.Range("E51").Value = '  -1.28%  '
